$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Objetivos:" value (B/C) changes to the professor's name string
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13: used to hold only the orphan professor-name value (B/C, no label).
# It becomes the "Programa resumido:" row with a "Semestral" value.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# Row 14: becomes the "Short syllabus:" label-only row (clear old paragraph value).
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# Row 15: becomes "Programa:" with a date value.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2016"
$ws.Rows(15).RowHeight = 120

# Row 16: becomes the "Syllabus:" label-only row (clear old long paragraph value).
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# Row 17: becomes "Avaliação:" label-only row with default (non-custom) height.
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows(17).AutoFit()

# Row 18: becomes "Método:" with the professor's name value again.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows(18).RowHeight = 60

# Row 19: label becomes "Critério:" (value text unchanged).
$ws.Range("A19").Value = "Critério:"

# Row 20: label becomes "Norma de recuperação:" (value text unchanged).
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: label becomes "Bibliografia:" (value text unchanged), taller row.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows(21).RowHeight = 120

# Row 22 (old Bibliografia text row) is removed entirely.
$ws.Rows(22).Delete()
